$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biomass reactions")
$ws.Range("A1").Value = "Test"
